$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the window back to the left (it had drifted to a second-monitor offset)
$win = $wb.Windows.Item(1)
$win.Left = -120

# Update the "Bilder" task detail text: Product::ImageId -> Product::ImagePath
$ws.Range("B3").Value = "Dateisystem für Bilder, Product::ImagePath, ProductForm & ProductOverview: Bild anzeigen"

# Remove the "Produktdatenbank mit Infos und Bildern füllen" detail from "Datenbank füllen"
$ws.Range("B11").ClearContents()

# Remove the "Lokaler Datencache" row entirely (row 12), keeping the blank-row cell formatting
$ws.Range("A12:H12").ClearContents()

# Drop the "Sonstiges" / assignee columns (G and H) entirely
$ws.Columns("G:H").Delete()

# The remaining bold-ish (but visually default) formatting on these task/detail
# cells is no longer needed once the "Sonstiges" column bookkeeping is gone
$ws.Range("A2:B2").ClearFormats()
$ws.Range("A4:B4").ClearFormats()
$ws.Range("A7:B7").ClearFormats()

# Reselect B9, matching where the cursor ended up after the reorganisation
$ws.Range("B9").Select()
